# Update chapter title numbering from "3.x"/"2.x" to "6.x" across slides.
$p = $ppt.ActivePresentation

# Slide 3: Title "3.1 Code" -> "6.1 Code"
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "6.1 Code"

# Slide 4: Title "2.1 Code" -> "6.1 Code"
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "6.1 Code"

# Slide 5: Title "3.2 Verify" -> "6.2 Verify"
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "6.2 Verify"

# Slide 6: Title "3.2 Verify" -> "6.2 Verify"
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "6.2 Verify"
